$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 373, shifting existing rows 373:396 down to 374:397.
$ws.Rows("373:373").Insert()

# Populate the newly inserted row 373 with the new weekly record.
$ws.Range("A373").Value = 4
$ws.Range("B373").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C373").Value = 'Los Lagos'
$ws.Range("D373").Value = 45265
$ws.Range("E373").Value = 10
$ws.Range("F373").Value = 'Fruta'
$ws.Range("G373").Value = 100109
$ws.Range("H373").Value = 'Uva'
$ws.Range("I373").Value = 100109001
$ws.Range("J373").Value = 'Uva'
$ws.Range("K373").Value = 'Superior Seedless'
$ws.Range("L373").Value = 'Primera'
$ws.Range("M373").Value = 150
$ws.Range("N373").Value = 23000
$ws.Range("O373").Value = 23000
$ws.Range("P373").Value = 23000
$ws.Range("Q373").Value = '$/bandeja 8 kilos'
$ws.Range("R373").Value = 'Provincia de Limarí'
$ws.Range("S373").Value = 2875
$ws.Range("T373").Value = 8
